# Update sample-suspension.xlsx:
#  - rename "preparation_temperature list" -> "preparation_condition list" with new values
#  - rename "storage_temperature list" -> "storage_method list" with new values
#  - update header text, comments and data validation on the "Export as TSV" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. preparation_temperature list -> preparation_condition list
# ---------------------------------------------------------------------------
$sPrepCond = $wb.Worksheets.Item("preparation_temperature list")
$sPrepCond.Name = "preparation_condition list"

$prepCondValues = @(
    "frozen in liquid nitrogen",
    "frozen in liquid nitrogen vapor",
    "frozen in ice",
    "frozen in dry ice",
    "frozen at -20 C",
    "ambient temperature",
    "unknown"
)
for ($i = 0; $i -lt $prepCondValues.Length; $i++) {
    $sPrepCond.Cells.Item($i + 1, 1).Value = $prepCondValues[$i]
}
$sPrepCond.Range("A8").ClearContents()

# ---------------------------------------------------------------------------
# 2. storage_temperature list -> storage_method list
# ---------------------------------------------------------------------------
$sStorMethod = $wb.Worksheets.Item("storage_temperature list")
$sStorMethod.Name = "storage_method list"

$storMethodValues = @(
    "frozen in liquid nitrogen",
    "frozen in liquid nitrogen vapor",
    "frozen in ice",
    "frozen in dry ice",
    "frozen at -80 C",
    "frozen at -20 C",
    "refrigerator",
    "ambient temperature",
    "incubated at 37 C",
    "none",
    "unknown"
)
for ($i = 0; $i -lt $storMethodValues.Length; $i++) {
    $sStorMethod.Cells.Item($i + 1, 1).Value = $storMethodValues[$i]
}
$sStorMethod.Range("A12").ClearContents()

# ---------------------------------------------------------------------------
# 3. "Export as TSV" sheet: headers, comments, data validation
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Export as TSV")

# -- Column G: preparation_condition --------------------------------------
$ws.Range("G1").Value = "preparation_condition"
$ws.Range("G1").Comment.Text("The condition under which the preparation occurred, such as whether the sample was placed in dry ice during the preparation.")

$gValidation = $ws.Range("G2:G1048576").Validation
$gValidation.Modify(3, 1, 1, "='preparation_condition list'!`$A`$1:`$A`$7")
$gValidation.ErrorTitle = "Value must come from list"
$gValidation.ErrorMessage = "Value must come from preparation_condition list."

# -- Column K: storage_method -----------------------------------------------
$ws.Range("K1").Value = "storage_method"
$ws.Range("K1").Comment.Text("The method by which the sample was stored, after preparation and before the assay was performed.")

$kValidation = $ws.Range("K2:K1048576").Validation
$kValidation.Modify(3, 1, 1, "='storage_method list'!`$A`$1:`$A`$11")
$kValidation.ErrorTitle = "Value must come from list"
$kValidation.ErrorMessage = "Value must come from storage_method list."
